$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Towns")

# Add new "order" column header in I1
$ws.Range("I1").Value = "order"

# Map of row -> order value (row 2..14 correspond to towns in sheet order)
$orderValues = @{
    2  = 1
    3  = 2
    4  = 13
    5  = 3
    6  = 8
    7  = 5
    8  = 6
    9  = 7
    10 = 4
    11 = 12
    12 = 9
    13 = 11
    14 = 10
}

foreach ($row in $orderValues.Keys) {
    $ws.Cells.Item($row, 9).Value = $orderValues[$row]
}

# Update selection to match target (K9) as seen in diff
$ws.Range("K9").Select()
